$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Registers sheet: move the addressing-mode registers (CREG, MREG, MREGC)
# down to the end of the table, shifting ENCREG1 / ENCREG2 / FLAGR / IHBASE
# up by two rows to fill the gap.
# ----------------------------------------------------------------------------
$regs = $wb.Worksheets.Item("Registers")

# Unmerge the two 2-row label/description blocks that are about to move.
$regs.Range("C10:C11").UnMerge()
$regs.Range("D10:D11").UnMerge()

# Row 8: was CREG -> becomes ENCREG1 (+ ENCREG12 note + description)
$regs.Range("B8").Value = "ENCREG1"
$regs.Range("C8").Value = "ENCREG12"
$regs.Range("D8").Value = "Will XOR next read OP byte with &0xFF, then shift with carry-around"

# Row 9: was MREG -> becomes ENCREG2
$regs.Range("B9").Value = "ENCREG2"
$regs.Range("C9").Value = $null
$regs.Range("D9").Value = $null

# Merge the moved ENCREG1/ENCREG12 label across its new home (rows 8-9).
$regs.Range("C8:C9").Merge()

# Apply "left + vertical-center" alignment to the cells that now carry the
# ENCREG1 description / ENCREG12 label (matches the other note cells).
$regs.Range("C8:C9").HorizontalAlignment = -4131
$regs.Range("C8:C9").VerticalAlignment = -4108
$regs.Range("D8:D9").HorizontalAlignment = -4131
$regs.Range("D8:D9").VerticalAlignment = -4108

# Row 10: was ENCREG1 -> becomes FLAGR
$regs.Range("B10").Value = "FLAGR"
$regs.Range("C10").Value = $null
$regs.Range("D10").Value = "Flag register"

# Row 11: was ENCREG2 -> becomes IHBASE
$regs.Range("B11").Value = "IHBASE"
$regs.Range("D11").Value = "Interrupt handler base (if 0, interrupts are off)"

# Row 12: was FLAGR -> now empty
$regs.Range("B12").Value = $null
$regs.Range("D12").Value = $null

# Row 13: was IHBASE -> now empty
$regs.Range("B13").Value = $null
$regs.Range("D13").Value = $null

# Row 14: was MREGC -> now empty
$regs.Range("B14").Value = $null
$regs.Range("D14").Value = $null

# Row 15: now CREG (was empty)
$regs.Range("B15").Value = "CREG"
$regs.Range("D15").Value = "Will swap to const mode in instructions"

# Row 16: now MREG (was empty)
$regs.Range("B16").Value = "MREG"
$regs.Range("D16").Value = "Will swap to RAM mode in instructions (address in reg)"

# Row 17: now MREGC (was empty)
$regs.Range("B17").Value = "MREGC"
$regs.Range("D17").Value = "Will swap to RAM mode in instructions (address const)"

# ----------------------------------------------------------------------------
# View state: Registers tab becomes the active sheet/tab, with a new
# selection; Opcodes loses its "active" state and gets a new selection too.
# ----------------------------------------------------------------------------
$opcodes = $wb.Worksheets.Item("Opcodes")
$opcodes.Activate()
$opcodes.Range("B58").Select()

$regs.Activate()
$regs.Range("C16").Select()
